$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '51.802.51'
$ws.Cells.Item(2, 5).Value = '  +0.57%  '

$ws.Cells.Item(3, 4).Value = '2.827.93'
$ws.Cells.Item(3, 5).Value = '  +1.81%  '

$ws.Cells.Item(4, 5).Value = '  -0.07%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '350.30'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.61%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '112.95'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +4.22%  '

$ws.Cells.Item(7, 5).Value = '  +1.36%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.619'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +3.78%  '

$ws.Cells.Item(10, 5).Value = '  +1.18%  '

$ws.Cells.Item(11, 5).Value = '  -0.76%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0849'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.67%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '20.05'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -0.15%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.77'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.54%  '

$ws.Cells.Item(15, 4).Value = '3.271.24'
$ws.Cells.Item(15, 5).Value = '  +1.72%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.980'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +6.14%  '

$ws.Cells.Item(17, 4).Value = '2.839.90'
$ws.Cells.Item(17, 5).Value = '  +2.30%  '

$ws.Cells.Item(18, 4).Value = '51.788.15'
$ws.Cells.Item(18, 5).Value = '  +0.56%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.43'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +11.04%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.61'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.61%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.41'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.18%  '

$ws.Cells.Item(22, 4).Value = '0.0₃0973'
$ws.Cells.Item(22, 5).Value = '  +1.10%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '70.42'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.77%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '268.99'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.10%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.75'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.72%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '26.26'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.71%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.00'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.01%  '

$ws.Cells.Item(28, 5).Value = '  +1.05%  '

$ws.Cells.Item(29, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '38.76'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +6.73%  '

$ws.Cells.Item(30, 2).Value = 'Cosmos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '10.53'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +3.01%  '

$ws.Cells.Item(31, 5).Value = '  +2.20%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.35'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +2.21%  '

$ws.Cells.Item(33, 5).Value = '  +1.68%  '

$ws.Cells.Item(34, 5).Value = '  +8.08%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '5.63'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.04%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0449'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.74%  '

$ws.Cells.Item(37, 5).Value = '  -0.07%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '18.91'
$ws.Cells.Item(38, 4).Style = "Normal"

$ws.Cells.Item(39, 5).Value = '  +2.20%  '

$ws.Cells.Item(40, 5).Value = '  +2.47%  '

$ws.Cells.Item(41, 5).Value = '  +1.35%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.52'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.21%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '122.86'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +2.14%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.23'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.92%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '22.08'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.09%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.55'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +9.65%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.50'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +7.66%  '

$ws.Cells.Item(48, 4).Value = '2.167.72'
$ws.Cells.Item(48, 5).Value = '  +2.89%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.248'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +21.86%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.949'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +5.05%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '5.53'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +2.05%  '
